$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at 130 (shifts existing rows 130-241 down to 131-242,
# matching the author's addition of a new weekly price observation).
$ws.Rows(130).Insert()

# Populate the newly inserted row 130 with the new record's data.
$ws.Cells.Item(130, 1).Value = 8
$ws.Cells.Item(130, 2).Value = "Terminal La Palmera de La Serena"
$ws.Cells.Item(130, 3).Value = "Coquimbo"
$ws.Cells.Item(130, 4).Value = 44827
$ws.Cells.Item(130, 5).Value = 4
$ws.Cells.Item(130, 6).Value = 100112037
$ws.Cells.Item(130, 7).Value = "Cebollín"
$ws.Cells.Item(130, 8).Value = "Sin especificar"
$ws.Cells.Item(130, 9).Value = "Primera"
$ws.Cells.Item(130, 10).Value = 2000
$ws.Cells.Item(130, 11).Value = 1400
$ws.Cells.Item(130, 12).Value = 1600
$ws.Cells.Item(130, 13).Value = 1500
$ws.Cells.Item(130, 14).Value = "`$/paquete 6 unidades"
$ws.Cells.Item(130, 15).Value = "Provincia del Elquí"
$ws.Cells.Item(130, 16).Value = 250
$ws.Cells.Item(130, 17).Value = 6
$ws.Cells.Item(130, 18).Value = "Hortaliza"
